# Daily attendance processing - 2025-10-30 16:26:49
# Normalizes the "Recorded By" (column G) values on the active worksheet:
#   - "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   - "system, backup@backdoor.com, System" -> "backup@backdoor.com, system, System"
# Other values (e.g. single names, or different combinations) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value2 = "backup@backdoor.com, system, System"
    }
}
